# Update FTSE 100 ticker list (Daten aktualisiert am 2023-09-17)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove Abrdn (ABDN) - no longer in the index
$ws.Rows("3:3").Delete()

# 2. Insert Dechra Pharmaceuticals (DPH) before Diageo (DGE), now at row 30
$ws.Rows("30:30").Insert()
$ws.Range("A30").Value = "DPH"
$ws.Range("B30").Value = "Dechra Pharmaceuticals"
$ws.Range("C30").Value = "Pharmaceuticals & Biotechnology"

# 3. Insert Diploma (DPLM) before Endeavour Mining (EDV), now at row 32
$ws.Rows("32:32").Insert()
$ws.Range("A32").Value = "DPLM"
$ws.Range("B32").Value = "Diploma"
$ws.Range("C32").Value = "Industrial Support Services"

# 4. Replace Hiscox (HSX) with Hikma Pharmaceuticals (HIK), now at row 45
$ws.Range("A45").Value = "HIK"
$ws.Range("B45").Value = "Hikma Pharmaceuticals"
$ws.Range("C45").Value = "Pharmaceuticals and Biotechnology"

# 5. Remove Johnson Matthey (JMAT) - no longer in the index, now at row 54
$ws.Rows("54:54").Delete()

# 6. Insert Marks & Spencer (MKS) before Melrose Industries (MRO), now at row 60
$ws.Rows("60:60").Insert()
$ws.Range("A60").Value = "MKS"
$ws.Range("B60").Value = "Marks & Spencer"
$ws.Range("C60").Value = "Diversified Retailers"
